$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.908.38"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.572.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.994"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.72"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.993"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.23"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +5.94%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.796.50"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.565.58"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.902.51"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.36"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.70"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +6.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0704"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.994"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.11"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.27"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.91%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.51"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.17"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.57"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0473"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.412.57"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.77%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.42"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.805"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.993"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.58"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.970"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.90%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.50%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.707.00"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.57"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.51%  "
